# edit.ps1 -- applies the "Tasks.docx" change described by the diff:
#   * "Website responsive machen" keeps its text but loses the bold
#     paragraph-mark formatting (becomes a plain paragraph).
#   * "Flaggen als Dropdown anbieten" is removed entirely.
#   * "transparentes HG Bild (auf Schramm warten für bild) " is shortened
#     to "transparentes HG Bild" (still a bold paragraph), and the
#     "_GoBack" bookmark (previously sitting between "durch" and
#     "newDate" further down the document) is moved to the end of this
#     paragraph.

$d = $word.ActiveDocument

$wdW = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function Find-ParagraphIndex($prefix) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text -like ($prefix + "*")) {
            return $i
        }
    }
    return $null
}

# 1) Move the "_GoBack" bookmark: delete it from its old location (it will
#    be re-created at the end of the "transparentes HG Bild" paragraph
#    below).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2) "Website responsive machen" paragraph: rewrite it without the bold
#    paragraph-mark formatting (<w:pPr><w:rPr><w:b/></w:rPr></w:pPr>),
#    keeping the same two runs (text + trailing space).
$iWebsite = Find-ParagraphIndex "Website responsive machen"
$pWebsite = $d.Paragraphs.Item($iWebsite)
$pWebsite.Range.InsertXML(
    '<w:p xmlns:w="' + $wdW + '">' +
        '<w:r><w:t>Website responsive machen</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '</w:p>'
)

# 3) Remove the "Flaggen als Dropdown anbieten" paragraph entirely.
$iFlaggen = Find-ParagraphIndex "Flaggen als Dropdown anbieten"
$pFlaggen = $d.Paragraphs.Item($iFlaggen)
$pFlaggen.Range.Delete()

# 4) Shorten "transparentes HG Bild (auf Schramm warten für bild) " down to
#    "transparentes HG Bild", keep the bold paragraph-mark formatting and
#    append the relocated "_GoBack" bookmark at the paragraph's end.
#    (Re-locate by text: the previous delete shifted paragraph indices.)
$iHgBild = Find-ParagraphIndex "transparentes HG Bild"
$pHgBild = $d.Paragraphs.Item($iHgBild)
$pHgBild.Range.InsertXML(
    '<w:p xmlns:w="' + $wdW + '">' +
        '<w:pPr><w:rPr><w:b/></w:rPr></w:pPr>' +
        '<w:r><w:t>transparentes HG Bild</w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
        '<w:bookmarkEnd w:id="0"/>' +
    '</w:p>'
)
